$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9910886883735657
$ws.Range("B1").Value = 1.911357402801514
$ws.Range("C1").Value = 5.38608455657959
$ws.Range("D1").Value = 2.288283348083496
$ws.Range("E1").Value = 1.302746415138245
